# Appends newly-logged sensor readings to the PIR, Humidity and Temperature
# sheets of the SeniorConnect master log (auto-update from the sensor feed).
#
# Dates (and the "NN.N%" humidity readings, which Excel's auto-detection
# also treats as numbers) are written with a leading apostrophe and then
# restyled back to "Normal" so Excel keeps them as plain text - e.g.
# "2026-02-01" / "80.1%" - instead of silently coercing them into a date
# serial / fractional number, matching how every other row in these sheets
# already stores its columns (plain text, no special number format).

$wb = $excel.ActiveWorkbook

# --- PIR sheet: rows 69-80 (Bathroom motion sensor) ---------------------
$pir = $wb.Worksheets.Item("PIR")
$pirData = @(
    @(69, "18:34:11", "No Motion",       "Inactive"),
    @(70, "18:34:13", "No Motion",       "Inactive"),
    @(71, "18:34:18", "No Motion",       "Inactive"),
    @(72, "18:34:23", "No Motion",       "Inactive"),
    @(73, "18:34:28", "No Motion",       "Inactive"),
    @(74, "18:34:33", "Motion Detected", "Active"),
    @(75, "18:34:41", "No Motion",       "Inactive"),
    @(76, "18:34:46", "No Motion",       "Inactive"),
    @(77, "18:34:51", "No Motion",       "Inactive"),
    @(78, "18:34:55", "Motion Detected", "Active"),
    @(79, "18:35:03", "No Motion",       "Inactive"),
    @(80, "18:35:08", "No Motion",       "Inactive")
)
foreach ($r in $pirData) {
    $row = $r[0]
    $pir.Cells.Item($row, 1).Value = "'2026-02-01"
    $pir.Cells.Item($row, 1).Style = "Normal"
    $pir.Cells.Item($row, 2).Value = $r[1]
    $pir.Cells.Item($row, 3).Value = "18:00"
    $pir.Cells.Item($row, 4).Value = "Bathroom"
    $pir.Cells.Item($row, 5).Value = $r[2]
    $pir.Cells.Item($row, 6).Value = $r[3]
}

# --- Humidity sheet: rows 128-138 (Bathroom humidity sensor) ------------
$humidity = $wb.Worksheets.Item("Humidity")
$humidityData = @(
    @(128, "18:34:14", "80.1%"),
    @(129, "18:34:19", "81.0%"),
    @(130, "18:34:24", "79.9%"),
    @(131, "18:34:34", "79.9%"),
    @(132, "18:34:39", "82.1%"),
    @(133, "18:34:44", "88.3%"),
    @(134, "18:34:49", "87.4%"),
    @(135, "18:34:54", "94.3%"),
    @(136, "18:34:59", "94.7%"),
    @(137, "18:35:04", "84.0%"),
    @(138, "18:35:09", "84.4%")
)
foreach ($r in $humidityData) {
    $row = $r[0]
    $humidity.Cells.Item($row, 1).Value = "'2026-02-01"
    $humidity.Cells.Item($row, 1).Style = "Normal"
    $humidity.Cells.Item($row, 2).Value = $r[1]
    $humidity.Cells.Item($row, 3).Value = "18:00"
    $humidity.Cells.Item($row, 4).Value = "Bathroom"
    $humidity.Cells.Item($row, 5).Value = "'" + $r[2]
    $humidity.Cells.Item($row, 5).Style = "Normal"
    $humidity.Cells.Item($row, 6).Value = "Active"
}

# --- Temperature sheet: rows 128-138 (Bathroom temperature sensor) ------
$temperature = $wb.Worksheets.Item("Temperature")
$temperatureData = @(
    @(128, "18:34:14", "29.1C"),
    @(129, "18:34:20", "29.1C"),
    @(130, "18:34:24", "29.1C"),
    @(131, "18:34:34", "29.1C"),
    @(132, "18:34:39", "29.1C"),
    @(133, "18:34:44", "29.2C"),
    @(134, "18:34:49", "29.2C"),
    @(135, "18:34:54", "29.3C"),
    @(136, "18:34:59", "29.3C"),
    @(137, "18:35:04", "29.3C"),
    @(138, "18:35:09", "29.3C")
)
foreach ($r in $temperatureData) {
    $row = $r[0]
    $temperature.Cells.Item($row, 1).Value = "'2026-02-01"
    $temperature.Cells.Item($row, 1).Style = "Normal"
    $temperature.Cells.Item($row, 2).Value = $r[1]
    $temperature.Cells.Item($row, 3).Value = "18:00"
    $temperature.Cells.Item($row, 4).Value = "Bathroom"
    $temperature.Cells.Item($row, 5).Value = $r[2]
    $temperature.Cells.Item($row, 6).Value = "Active"
}
